# Generate Report for Handoff
# For the zh-cn and de-de localization-status sheets, rows 4-7 (the four
# files that just got their handoff xliff regenerated) move from
# Priority "low" to "ht", and their "Latest Handoff Datetime" is bumped
# to the new handoff generation time. The Overview sheet's "Latest HO
# Xliff Generate Date" column shares that same string, so it picks up
# the newer (de-de) timestamp too.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
for ($row = 4; $row -le 7; $row++) {
    $ws_zhcn.Cells.Item($row, 5).Value = "ht"
    $ws_zhcn.Cells.Item($row, 8).Value = "2016-11-09 07:34:10"
}

$ws_dede = $wb.Worksheets.Item("de-de")
for ($row = 4; $row -le 7; $row++) {
    $ws_dede.Cells.Item($row, 5).Value = "ht"
    $ws_dede.Cells.Item($row, 8).Value = "2016-11-09 07:34:25"
}

$ws_overview = $wb.Worksheets.Item("Overview")
for ($row = 4; $row -le 7; $row++) {
    $ws_overview.Cells.Item($row, 7).Value = "2016-11-09 07:34:25"
}
